$d = $word.ActiveDocument

$d.Content.Find.Execute("Brute Force", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Breadth First Search and Depth First Search", 2)
